# Update computed market-price / profit columns (H:N) on the per-job Leve
# profit sheets. Values reflect refreshed Universalis price snapshots;
# row/item identity confirmed via the "Leve Item ID" (G) column.
$wb = $excel.ActiveWorkbook

# ===================== ALC =====================
$ws = $wb.Worksheets.Item("ALC")
# Row 33: Glazed and Confused
$ws.Range("H33").Value = 651.88
$ws.Range("I33").Value = 507.25
$ws.Range("K33").Value = 507.25
$ws.Range("M33").Value = -278.25
# Row 40: Stuck in the Moment
$ws.Range("H40").Value = 2011.1111
$ws.Range("I40").Value = 2000
$ws.Range("J40").Value = 2033.3334
$ws.Range("K40").Value = 2000
$ws.Range("L40").Value = 2033.3334
$ws.Range("M40").Value = -1825
$ws.Range("N40").Value = -2383.3334
# Row 74: Adhesive of Antipathy
$ws.Range("H74").Value = 4592
$ws.Range("I74").Value = 3986
$ws.Range("J74").Value = 4895
$ws.Range("K74").Value = 3986
$ws.Range("L74").Value = 4895
$ws.Range("M74").Value = -3050
$ws.Range("N74").Value = -6767
# Row 76: Warding Off Temptation
$ws.Range("H76").Value = 2964.3076
$ws.Range("I76").Value = 2893.1
$ws.Range("J76").Value = 3201.6667
$ws.Range("K76").Value = 2893.1
$ws.Range("L76").Value = 3201.6667
$ws.Range("M76").Value = -2578.1
$ws.Range("N76").Value = -3831.6667
# Row 77: It's Gonna Grow Back (L)
$ws.Range("H77").Value = 4592
$ws.Range("I77").Value = 3986
$ws.Range("J77").Value = 4895
$ws.Range("K77").Value = 19930
$ws.Range("L77").Value = 24475
$ws.Range("M77").Value = -15250
$ws.Range("N77").Value = -33835
# Row 79: The Garden of Arcane Delights (L)
$ws.Range("H79").Value = 2964.3076
$ws.Range("I79").Value = 2893.1
$ws.Range("J79").Value = 3201.6667
$ws.Range("K79").Value = 2893.1
$ws.Range("L79").Value = 3201.6667
$ws.Range("M79").Value = -1801.1
$ws.Range("N79").Value = -5385.6667
# Row 134: Binding Spells
$ws.Range("H134").Value = 21196.363
$ws.Range("J134").Value = 21196.363
$ws.Range("L134").Value = 21196.363
$ws.Range("N134").Value = -31336.363
# Row 141: Remedy for Reason
$ws.Range("H141").Value = 383842.62
$ws.Range("I141").Value = 1408.0834
$ws.Range("J141").Value = 560350.9
$ws.Range("K141").Value = 4224.2502
$ws.Range("L141").Value = 1681052.7
$ws.Range("M141").Value = 955.7497999999996
$ws.Range("N141").Value = -1691412.7

# ===================== ARM =====================
$ws = $wb.Worksheets.Item("ARM")
# Row 9: Headbangers' Thrall
$ws.Range("H9").Value = 19300
$ws.Range("J9").Value = 19300
$ws.Range("L9").Value = 19300
$ws.Range("N9").Value = -19640
# Row 20: Cover Girl
$ws.Range("H20").Value = 19300
$ws.Range("J20").Value = 19300
$ws.Range("L20").Value = 19300
$ws.Range("N20").Value = -19840
# Row 23: A Well-rounded Crew
$ws.Range("H23").Value = 62754
$ws.Range("J23").Value = 52402.8
$ws.Range("L23").Value = 52402.8
$ws.Range("N23").Value = -52920.8
# Row 37: Get Shirty
$ws.Range("H37").Value = 16350
$ws.Range("J37").Value = 16350
$ws.Range("L37").Value = 16350
$ws.Range("N37").Value = -16896
# Row 44: Very Slow Array
$ws.Range("H44").Value = 26049
$ws.Range("I44").Value = 0
$ws.Range("J44").Value = 26049
$ws.Range("K44").Value = 0
$ws.Range("L44").Value = 26049
$ws.Range("M44").ClearContents()
$ws.Range("N44").Value = -27025
# Row 55: Employee Retention
$ws.Range("H55").Value = 19384.666
# Row 57: Cobalt Aforethought
$ws.Range("H57").Value = 10000
$ws.Range("I57").Value = 10000
$ws.Range("K57").Value = 10000
$ws.Range("M57").Value = -9516
# Row 61: Dealing with the Tough Stuff
$ws.Range("H61").Value = 2474.2083
$ws.Range("I61").Value = 1962.5
$ws.Range("J61").Value = 4009.3333
$ws.Range("K61").Value = 1962.5
$ws.Range("L61").Value = 4009.3333
$ws.Range("M61").Value = -1750.5
$ws.Range("N61").Value = -4433.3333
# Row 80: A Squire to Inspire
$ws.Range("H80").Value = 26764.75
$ws.Range("J80").Value = 26764.75
$ws.Range("L80").Value = 26764.75
$ws.Range("N80").Value = -28760.75
# Row 83: All's Fair in Highborn Assassination (L)
$ws.Range("H83").Value = 26764.75
$ws.Range("J83").Value = 26764.75
$ws.Range("L83").Value = 80294.25
$ws.Range("N83").Value = -90278.25
# Row 132: Don't Bore Me, Ore Me
$ws.Range("H132").Value = 2222.55
$ws.Range("I132").Value = 1616.5333
$ws.Range("J132").Value = 4040.6
$ws.Range("K132").Value = 4849.5999
$ws.Range("L132").Value = 12121.8
$ws.Range("M132").Value = -2319.5999
$ws.Range("N132").Value = -17181.8
# Row 136: Metal with Mettle
$ws.Range("H136").Value = 2474.2083
$ws.Range("I136").Value = 1962.5
$ws.Range("J136").Value = 4009.3333
$ws.Range("K136").Value = 5887.5
$ws.Range("L136").Value = 12027.9999
$ws.Range("M136").Value = -3337.5
$ws.Range("N136").Value = -17127.9999

# ===================== BSM =====================
$ws = $wb.Worksheets.Item("BSM")
# Row 22: Riveting Run
$ws.Range("H22").Value = 300
$ws.Range("I22").Value = 300
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 300
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -127
$ws.Range("N22").ClearContents()

# ===================== CUL =====================
$ws = $wb.Worksheets.Item("CUL")
# Row 69: Loving That Muffin Top
$ws.Range("H69").Value = 184168.67
$ws.Range("I69").Value = 670.6667
$ws.Range("K69").Value = 2012.0001
$ws.Range("M69").Value = -1201.0001
# Row 72: Muffin of the Morn (L)
$ws.Range("H72").Value = 184168.67
$ws.Range("I72").Value = 670.6667
$ws.Range("K72").Value = 6036.0003
$ws.Range("M72").Value = -1980.0003
# Row 74: The Nutcracker's Sweets
$ws.Range("H74").Value = 22200
$ws.Range("I74").Value = 1500
$ws.Range("J74").Value = 36000
$ws.Range("K74").Value = 4500
$ws.Range("L74").Value = 108000
$ws.Range("M74").Value = -3439
$ws.Range("N74").Value = -110122
# Row 75: Breakfast of Champions
$ws.Range("H75").Value = 2012.5
$ws.Range("I75").Value = 2012.5
$ws.Range("K75").Value = 6037.5
$ws.Range("M75").Value = -5039.5
# Row 76: Old Victories, New Tastes
$ws.Range("H76").Value = 3800
# Row 77: Time for a Midnight Snack (L)
$ws.Range("H77").Value = 22200
$ws.Range("I77").Value = 1500
$ws.Range("J77").Value = 36000
$ws.Range("K77").Value = 13500
$ws.Range("L77").Value = 324000
$ws.Range("M77").Value = -8196
$ws.Range("N77").Value = -334608
# Row 78: Emerald Soup for the Soul (L)
$ws.Range("H78").Value = 2012.5
$ws.Range("I78").Value = 2012.5
$ws.Range("K78").Value = 18112.5
$ws.Range("M78").Value = -13120.5
# Row 79: The Eats of Authenticity (L)
$ws.Range("H79").Value = 3800
# Row 80: Saucy for a Suitor
$ws.Range("H80").Value = 2880
$ws.Range("J80").Value = 2975
$ws.Range("L80").Value = 8925
$ws.Range("N80").Value = -10797
# Row 81: It Goes Down Smoothly
$ws.Range("H81").Value = 26600
$ws.Range("J81").Value = 26600
$ws.Range("L81").Value = 79800
$ws.Range("N81").Value = -82046
# Row 82: Persuasion of a Higher Power
$ws.Range("H82").Value = 3600
$ws.Range("I82").Value = 2000
$ws.Range("K82").Value = 6000
$ws.Range("M82").Value = -5594
# Row 83: Saved by the Sauce (L)
$ws.Range("H83").Value = 2880
$ws.Range("J83").Value = 2975
$ws.Range("L83").Value = 26775
$ws.Range("N83").Value = -36135
# Row 84: Quenching the Flame (L)
$ws.Range("H84").Value = 26600
$ws.Range("J84").Value = 26600
$ws.Range("L84").Value = 239400
$ws.Range("N84").Value = -250632
# Row 85: Loaves and Fishes (L)
$ws.Range("H85").Value = 3600
$ws.Range("I85").Value = 2000
$ws.Range("K85").Value = 6000
$ws.Range("M85").Value = -4596
# Row 86: Let's Not Get Sappy
$ws.Range("H86").Value = 675
$ws.Range("I86").Value = 350
$ws.Range("J86").Value = 1000
$ws.Range("K86").Value = 1050
$ws.Range("L86").Value = 3000
$ws.Range("M86").Value = 136
$ws.Range("N86").Value = -5372
# Row 87: Soup That Eats Like a Knight
$ws.Range("H87").Value = 13702
$ws.Range("I87").Value = 7950
$ws.Range("K87").Value = 23850
$ws.Range("M87").Value = -22602
# Row 89: Luxury Spillover (L)
$ws.Range("H89").Value = 675
$ws.Range("I89").Value = 350
$ws.Range("J89").Value = 1000
$ws.Range("K89").Value = 3150
$ws.Range("L89").Value = 9000
$ws.Range("M89").Value = 2778
$ws.Range("N89").Value = -20856
# Row 90: Like Ma Used to Make (L)
$ws.Range("H90").Value = 13702
$ws.Range("I90").Value = 7950
$ws.Range("K90").Value = 71550
$ws.Range("M90").Value = -65310
# Row 92: Oh No Udon
$ws.Range("H92").Value = 1163.5883
$ws.Range("I92").Value = 619.5
$ws.Range("J92").Value = 1331
$ws.Range("K92").Value = 1858.5
$ws.Range("L92").Value = 3993
$ws.Range("M92").Value = -610.5
$ws.Range("N92").Value = -6489
# Row 98: Sweet Kiss of Death
$ws.Range("H98").Value = 285.22223
$ws.Range("I98").Value = 200
$ws.Range("J98").Value = 455.66666
$ws.Range("K98").Value = 600
$ws.Range("L98").Value = 1366.99998
$ws.Range("M98").Value = 898
$ws.Range("N98").Value = -4362.999980000001
# Row 109: Cure for What Ails
$ws.Range("H109").Value = 682
$ws.Range("I109").Value = 205.71428
$ws.Range("J109").Value = 1793.3334
$ws.Range("K109").Value = 617.14284
$ws.Range("L109").Value = 5380.0002
$ws.Range("M109").Value = 422.85716
$ws.Range("N109").Value = -7460.0002
# Row 121: A Cookie for Your Troubles
$ws.Range("H121").Value = 8455.59
$ws.Range("I121").Value = 225.94444
$ws.Range("J121").Value = 15509.571
$ws.Range("K121").Value = 677.83332
$ws.Range("L121").Value = 46528.713
$ws.Range("M121").Value = 632.16668
$ws.Range("N121").Value = -49148.713

# ===================== GSM =====================
$ws = $wb.Worksheets.Item("GSM")
# Row 70: Sky Is the Limit
$ws.Range("H70").Value = 4418.385
$ws.Range("I70").Value = 4420
$ws.Range("J70").Value = 4416.1816
$ws.Range("K70").Value = 4420
$ws.Range("L70").Value = 4416.1816
$ws.Range("M70").Value = -4150
$ws.Range("N70").Value = -4956.1816
# Row 73: Hulls of Broken Dreams (L)
$ws.Range("H73").Value = 4418.385
$ws.Range("I73").Value = 4420
$ws.Range("J73").Value = 4416.1816
$ws.Range("K73").Value = 4420
$ws.Range("L73").Value = 4416.1816
$ws.Range("M73").Value = -3484
$ws.Range("N73").Value = -6288.1816
# Row 80: Needs More Prayerbell
$ws.Range("H80").Value = 4510.4
$ws.Range("I80").Value = 4048.6667
$ws.Range("J80").Value = 5203
$ws.Range("K80").Value = 4048.6667
$ws.Range("L80").Value = 5203
$ws.Range("M80").Value = -3050.6667
$ws.Range("N80").Value = -7199
# Row 83: With a Noise That Reaches Heaven (L)
$ws.Range("H83").Value = 4510.4
$ws.Range("I83").Value = 4048.6667
$ws.Range("J83").Value = 5203
$ws.Range("K83").Value = 20243.3335
$ws.Range("L83").Value = 26015
$ws.Range("M83").Value = -15251.3335
$ws.Range("N83").Value = -35999

# ===================== LTW =====================
$ws = $wb.Worksheets.Item("LTW")
# Row 35: No Risk, No Reward
$ws.Range("H35").Value = 10153.143
$ws.Range("I35").Value = 268
$ws.Range("K35").Value = 268
$ws.Range("M35").Value = 68
# Row 58: Handle with Care
$ws.Range("H58").Value = 2721.25
$ws.Range("I58").Value = 1965
$ws.Range("J58").Value = 4990
$ws.Range("K58").Value = 1965
$ws.Range("L58").Value = 4990
$ws.Range("M58").Value = -1705
$ws.Range("N58").Value = -5510

